# Apply updated currentAveragePrice / profit figures per the scheduled market-data refresh.
# Values correspond to columns H:N (currentAveragePrice.. LeveProfitHQ) on each sheet's Leve table.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 227.6  # H9
$ws.Cells.Item(9, 9).Value = 190.67857  # I9
$ws.Cells.Item(9, 10).Value = 375.2857  # J9
$ws.Cells.Item(9, 11).Value = 190.67857  # K9
$ws.Cells.Item(9, 12).Value = 375.2857  # L9
$ws.Cells.Item(9, 13).Value = -21.67857000000001  # M9
$ws.Cells.Item(9, 14).Value = -713.2857  # N9
$ws.Cells.Item(11, 8).Value = 477.45456  # H11
$ws.Cells.Item(11, 9).Value = 477.45456  # I11
$ws.Cells.Item(11, 11).Value = 477.45456  # K11
$ws.Cells.Item(11, 13).Value = -337.45456  # M11
$ws.Cells.Item(12, 8).Value = 1095.7142  # H12
$ws.Cells.Item(12, 9).Value = 833.6  # I12
$ws.Cells.Item(12, 11).Value = 833.6  # K12
$ws.Cells.Item(12, 13).Value = -663.6  # M12
$ws.Cells.Item(19, 8).Value = 17289  # H19
$ws.Cells.Item(19, 9).Value = 996.3333  # I19
$ws.Cells.Item(19, 11).Value = 996.3333  # K19
$ws.Cells.Item(19, 13).Value = -821.3333  # M19
$ws.Cells.Item(40, 8).Value = 4010.8333  # H40
$ws.Cells.Item(40, 9).Value = 3474.75  # I40
$ws.Cells.Item(40, 10).Value = 5083  # J40
$ws.Cells.Item(40, 11).Value = 3474.75  # K40
$ws.Cells.Item(40, 12).Value = 5083  # L40
$ws.Cells.Item(40, 13).Value = -3299.75  # M40
$ws.Cells.Item(40, 14).Value = -5433  # N40
$ws.Cells.Item(58, 8).Value = 4354.45  # H58
$ws.Cells.Item(58, 9).Value = 63  # I58
$ws.Cells.Item(58, 10).Value = 7865.636  # J58
$ws.Cells.Item(58, 11).Value = 189  # K58
$ws.Cells.Item(58, 12).Value = 23596.908  # L58
$ws.Cells.Item(58, 13).Value = -39  # M58
$ws.Cells.Item(58, 14).Value = -23896.908  # N58
$ws.Cells.Item(64, 8).Value = 3813.6  # H64
$ws.Cells.Item(64, 9).Value = 3333.6667  # I64
$ws.Cells.Item(64, 10).Value = 4133.5557  # J64
$ws.Cells.Item(64, 11).Value = 3333.6667  # K64
$ws.Cells.Item(64, 12).Value = 4133.5557  # L64
$ws.Cells.Item(64, 13).Value = -3085.6667  # M64
$ws.Cells.Item(64, 14).Value = -4629.5557  # N64
$ws.Cells.Item(67, 8).Value = 3813.6  # H67
$ws.Cells.Item(67, 9).Value = 3333.6667  # I67
$ws.Cells.Item(67, 10).Value = 4133.5557  # J67
$ws.Cells.Item(67, 11).Value = 3333.6667  # K67
$ws.Cells.Item(67, 12).Value = 4133.5557  # L67
$ws.Cells.Item(67, 13).Value = -2475.6667  # M67
$ws.Cells.Item(67, 14).Value = -5849.5557  # N67
$ws.Cells.Item(86, 8).Value = 3212.077  # H86
$ws.Cells.Item(86, 9).Value = 3306.2222  # I86
$ws.Cells.Item(86, 10).Value = 3000.25  # J86
$ws.Cells.Item(86, 11).Value = 3306.2222  # K86
$ws.Cells.Item(86, 12).Value = 3000.25  # L86
$ws.Cells.Item(86, 13).Value = -2183.2222  # M86
$ws.Cells.Item(86, 14).Value = -5246.25  # N86
$ws.Cells.Item(89, 8).Value = 3212.077  # H89
$ws.Cells.Item(89, 9).Value = 3306.2222  # I89
$ws.Cells.Item(89, 10).Value = 3000.25  # J89
$ws.Cells.Item(89, 11).Value = 16531.111  # K89
$ws.Cells.Item(89, 12).Value = 15001.25  # L89
$ws.Cells.Item(89, 13).Value = -10915.111  # M89
$ws.Cells.Item(89, 14).Value = -26233.25  # N89
$ws.Cells.Item(132, 8).Value = 5331.244  # H132
$ws.Cells.Item(132, 9).Value = 3378.56  # I132
$ws.Cells.Item(132, 11).Value = 10135.68  # K132
$ws.Cells.Item(132, 13).Value = -7605.68  # M132
$ws.Cells.Item(135, 8).Value = 2456.44  # H135
$ws.Cells.Item(135, 9).Value = 1793.8422  # I135
$ws.Cells.Item(135, 11).Value = 16144.5798  # K135
$ws.Cells.Item(135, 13).Value = -13609.5798  # M135
$ws.Cells.Item(137, 8).Value = 3006.1936  # H137
$ws.Cells.Item(137, 9).Value = 2963.25  # I137
$ws.Cells.Item(137, 10).Value = 3084.2727  # J137
$ws.Cells.Item(137, 11).Value = 8889.75  # K137
$ws.Cells.Item(137, 12).Value = 9252.8181  # L137
$ws.Cells.Item(137, 13).Value = -6339.75  # M137
$ws.Cells.Item(137, 14).Value = -14352.8181  # N137
$ws.Cells.Item(138, 8).Value = 3149.6304  # H138
$ws.Cells.Item(138, 9).Value = 2882.625  # I138
$ws.Cells.Item(138, 10).Value = 3292.0334  # J138
$ws.Cells.Item(138, 11).Value = 8647.875  # K138
$ws.Cells.Item(138, 12).Value = 9876.100199999999  # L138
$ws.Cells.Item(138, 13).Value = -3507.875  # M138
$ws.Cells.Item(138, 14).Value = -20156.1002  # N138
$ws.Cells.Item(141, 8).Value = 15121.5  # H141
$ws.Cells.Item(141, 9).Value = 12985.3  # I141
$ws.Cells.Item(141, 11).Value = 38955.89999999999  # K141
$ws.Cells.Item(141, 13).Value = -33775.89999999999  # M141

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3817.9365  # H32
$ws.Cells.Item(32, 9).Value = 2768.638  # I32
$ws.Cells.Item(32, 11).Value = 2768.638  # K32
$ws.Cells.Item(32, 13).Value = -2481.638  # M32
$ws.Cells.Item(61, 8).Value = 2706.077  # H61
$ws.Cells.Item(61, 9).Value = 1585.5  # I61
$ws.Cells.Item(61, 10).Value = 7828.7144  # J61
$ws.Cells.Item(61, 11).Value = 1585.5  # K61
$ws.Cells.Item(61, 12).Value = 7828.7144  # L61
$ws.Cells.Item(61, 13).Value = -1373.5  # M61
$ws.Cells.Item(61, 14).Value = -8252.714400000001  # N61
$ws.Cells.Item(64, 8).Value = 49599.4  # H64
$ws.Cells.Item(64, 10).Value = 49599.4  # J64
$ws.Cells.Item(64, 12).Value = 49599.4  # L64
$ws.Cells.Item(64, 14).Value = -50095.4  # N64
$ws.Cells.Item(67, 8).Value = 49599.4  # H67
$ws.Cells.Item(67, 10).Value = 49599.4  # J67
$ws.Cells.Item(67, 12).Value = 49599.4  # L67
$ws.Cells.Item(67, 14).Value = -51315.4  # N67
$ws.Cells.Item(74, 8).Value = 7166.6665  # H74
$ws.Cells.Item(74, 9).Value = 1500  # I74
$ws.Cells.Item(74, 10).Value = 10000  # J74
$ws.Cells.Item(74, 11).Value = 1500  # K74
$ws.Cells.Item(74, 12).Value = 10000  # L74
$ws.Cells.Item(74, 13).Value = -626  # M74
$ws.Cells.Item(74, 14).Value = -11748  # N74
$ws.Cells.Item(77, 8).Value = 7166.6665  # H77
$ws.Cells.Item(77, 9).Value = 1500  # I77
$ws.Cells.Item(77, 10).Value = 10000  # J77
$ws.Cells.Item(77, 11).Value = 7500  # K77
$ws.Cells.Item(77, 12).Value = 50000  # L77
$ws.Cells.Item(77, 13).Value = -3132  # M77
$ws.Cells.Item(77, 14).Value = -58736  # N77
$ws.Cells.Item(92, 8).Value = 151333  # H92
$ws.Cells.Item(92, 10).Value = 151333  # J92
$ws.Cells.Item(92, 12).Value = 151333  # L92
$ws.Cells.Item(92, 14).Value = -156325  # N92
$ws.Cells.Item(97, 8).Value = 4826.5654  # H97
$ws.Cells.Item(97, 9).Value = 5253.8237  # I97
$ws.Cells.Item(97, 11).Value = 5253.8237  # K97
$ws.Cells.Item(97, 13).Value = -4757.8237  # M97
$ws.Cells.Item(102, 8).Value = 5624.2  # H102
$ws.Cells.Item(102, 9).Value = 4598.909  # I102
$ws.Cells.Item(102, 11).Value = 4598.909  # K102
$ws.Cells.Item(102, 13).Value = -2976.909  # M102
$ws.Cells.Item(122, 8).Value = 9459.729499999999  # H122
$ws.Cells.Item(122, 9).Value = 8838.857  # I122
$ws.Cells.Item(122, 11).Value = 26516.571  # K122
$ws.Cells.Item(122, 13).Value = -24066.571  # M122
$ws.Cells.Item(132, 8).Value = 1613.8971  # H132
$ws.Cells.Item(132, 9).Value = 1314.0339  # I132
$ws.Cells.Item(132, 10).Value = 3579.6667  # J132
$ws.Cells.Item(132, 11).Value = 3942.1017  # K132
$ws.Cells.Item(132, 12).Value = 10739.0001  # L132
$ws.Cells.Item(132, 13).Value = -1412.1017  # M132
$ws.Cells.Item(132, 14).Value = -15799.0001  # N132
$ws.Cells.Item(136, 8).Value = 2706.077  # H136
$ws.Cells.Item(136, 9).Value = 1585.5  # I136
$ws.Cells.Item(136, 10).Value = 7828.7144  # J136
$ws.Cells.Item(136, 11).Value = 4756.5  # K136
$ws.Cells.Item(136, 12).Value = 23486.1432  # L136
$ws.Cells.Item(136, 13).Value = -2206.5  # M136
$ws.Cells.Item(136, 14).Value = -28586.1432  # N136

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1735.1666  # H86
$ws.Cells.Item(86, 9).Value = 1776.5834  # I86
$ws.Cells.Item(86, 10).Value = 1652.3334  # J86
$ws.Cells.Item(86, 11).Value = 1776.5834  # K86
$ws.Cells.Item(86, 12).Value = 1652.3334  # L86
$ws.Cells.Item(86, 13).Value = -653.5834  # M86
$ws.Cells.Item(86, 14).Value = -3898.3334  # N86
$ws.Cells.Item(89, 8).Value = 1735.1666  # H89
$ws.Cells.Item(89, 9).Value = 1776.5834  # I89
$ws.Cells.Item(89, 10).Value = 1652.3334  # J89
$ws.Cells.Item(89, 11).Value = 8882.916999999999  # K89
$ws.Cells.Item(89, 12).Value = 8261.666999999999  # L89
$ws.Cells.Item(89, 13).Value = -3266.916999999999  # M89
$ws.Cells.Item(89, 14).Value = -19493.667  # N89
$ws.Cells.Item(95, 8).Value = 46859.855  # H95
$ws.Cells.Item(95, 10).Value = 46859.855  # J95
$ws.Cells.Item(95, 12).Value = 46859.855  # L95
$ws.Cells.Item(95, 14).Value = -52351.855  # N95
$ws.Cells.Item(105, 8).Value = 2774.8  # H105
$ws.Cells.Item(105, 9).Value = 1779.1852  # I105
$ws.Cells.Item(105, 11).Value = 1779.1852  # K105
$ws.Cells.Item(105, 13).Value = -32.1851999999999  # M105
$ws.Cells.Item(124, 8).Value = 66660.664  # H124
$ws.Cells.Item(124, 10).Value = 66660.664  # J124
$ws.Cells.Item(124, 12).Value = 66660.664  # L124
$ws.Cells.Item(124, 14).Value = -76480.664  # N124
$ws.Cells.Item(126, 8).Value = 69999  # H126
$ws.Cells.Item(126, 10).Value = 69999  # J126
$ws.Cells.Item(126, 12).Value = 69999  # L126
$ws.Cells.Item(126, 14).Value = -79879  # N126
$ws.Cells.Item(134, 8).Value = 3233.3333  # H134
$ws.Cells.Item(134, 9).Value = 2597.0908  # I134
$ws.Cells.Item(134, 10).Value = 4233.143  # J134
$ws.Cells.Item(134, 11).Value = 7791.2724  # K134
$ws.Cells.Item(134, 12).Value = 12699.429  # L134
$ws.Cells.Item(134, 13).Value = -5256.2724  # M134
$ws.Cells.Item(134, 14).Value = -17769.429  # N134

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4556.5557  # H16
$ws.Cells.Item(16, 9).Value = 2002.75  # I16
$ws.Cells.Item(16, 10).Value = 6599.6  # J16
$ws.Cells.Item(16, 11).Value = 2002.75  # K16
$ws.Cells.Item(16, 12).Value = 6599.6  # L16
$ws.Cells.Item(16, 13).Value = -1715.75  # M16
$ws.Cells.Item(16, 14).Value = -7173.6  # N16
$ws.Cells.Item(31, 8).Value = 2662  # H31
$ws.Cells.Item(31, 9).Value = 1883.6666  # I31
$ws.Cells.Item(31, 11).Value = 1883.6666  # K31
$ws.Cells.Item(31, 13).Value = -1588.6666  # M31
$ws.Cells.Item(34, 8).Value = 2662  # H34
$ws.Cells.Item(34, 9).Value = 1883.6666  # I34
$ws.Cells.Item(34, 11).Value = 1883.6666  # K34
$ws.Cells.Item(34, 13).Value = -1681.6666  # M34
$ws.Cells.Item(58, 8).Value = 5737.8213  # H58
$ws.Cells.Item(58, 9).Value = 5315  # I58
$ws.Cells.Item(58, 11).Value = 5315  # K58
$ws.Cells.Item(58, 13).Value = -5112  # M58
$ws.Cells.Item(74, 8).Value = 42947.08  # H74
$ws.Cells.Item(74, 10).Value = 42947.08  # J74
$ws.Cells.Item(74, 12).Value = 42947.08  # L74
$ws.Cells.Item(74, 14).Value = -44695.08  # N74
$ws.Cells.Item(77, 8).Value = 42947.08  # H77
$ws.Cells.Item(77, 10).Value = 42947.08  # J77
$ws.Cells.Item(77, 12).Value = 128841.24  # L77
$ws.Cells.Item(77, 14).Value = -137577.24  # N77
$ws.Cells.Item(86, 8).Value = 7412462.5  # H86
$ws.Cells.Item(86, 9).Value = 9528166  # I86
$ws.Cells.Item(86, 10).Value = 7499.5  # J86
$ws.Cells.Item(86, 11).Value = 9528166  # K86
$ws.Cells.Item(86, 12).Value = 7499.5  # L86
$ws.Cells.Item(86, 13).Value = -9527043  # M86
$ws.Cells.Item(86, 14).Value = -9745.5  # N86
$ws.Cells.Item(89, 8).Value = 7412462.5  # H89
$ws.Cells.Item(89, 9).Value = 9528166  # I89
$ws.Cells.Item(89, 10).Value = 7499.5  # J89
$ws.Cells.Item(89, 11).Value = 47640830  # K89
$ws.Cells.Item(89, 12).Value = 37497.5  # L89
$ws.Cells.Item(89, 13).Value = -47635214  # M89
$ws.Cells.Item(89, 14).Value = -48729.5  # N89
$ws.Cells.Item(107, 8).Value = 2571.647  # H107
$ws.Cells.Item(107, 9).Value = 1434.9166  # I107
$ws.Cells.Item(107, 11).Value = 1434.9166  # K107
$ws.Cells.Item(107, 13).Value = 485.0834  # M107
$ws.Cells.Item(113, 8).Value = 4556.5557  # H113
$ws.Cells.Item(113, 9).Value = 2002.75  # I113
$ws.Cells.Item(113, 10).Value = 6599.6  # J113
$ws.Cells.Item(113, 11).Value = 2002.75  # K113
$ws.Cells.Item(113, 12).Value = 6599.6  # L113
$ws.Cells.Item(113, 13).Value = 167.25  # M113
$ws.Cells.Item(113, 14).Value = -10939.6  # N113
$ws.Cells.Item(122, 8).Value = 1888.9032  # H122
$ws.Cells.Item(122, 9).Value = 1465.5186  # I122
$ws.Cells.Item(122, 10).Value = 4746.75  # J122
$ws.Cells.Item(122, 11).Value = 4396.5558  # K122
$ws.Cells.Item(122, 12).Value = 14240.25  # L122
$ws.Cells.Item(122, 13).Value = -1946.5558  # M122
$ws.Cells.Item(122, 14).Value = -19140.25  # N122
$ws.Cells.Item(132, 8).Value = 7042.643  # H132
$ws.Cells.Item(132, 9).Value = 9051.684999999999  # I132
$ws.Cells.Item(132, 11).Value = 27155.055  # K132
$ws.Cells.Item(132, 13).Value = -24625.055  # M132
$ws.Cells.Item(134, 8).Value = 3465.4656  # H134
$ws.Cells.Item(134, 9).Value = 2648.3403  # I134
$ws.Cells.Item(134, 10).Value = 6956.8184  # J134
$ws.Cells.Item(134, 11).Value = 7945.0209  # K134
$ws.Cells.Item(134, 12).Value = 20870.4552  # L134
$ws.Cells.Item(134, 13).Value = -5410.0209  # M134
$ws.Cells.Item(134, 14).Value = -25940.4552  # N134
$ws.Cells.Item(136, 8).Value = 5737.8213  # H136
$ws.Cells.Item(136, 9).Value = 5315  # I136
$ws.Cells.Item(136, 11).Value = 15945  # K136
$ws.Cells.Item(136, 13).Value = -13395  # M136

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(88, 8).Value = 9117.777  # H88
$ws.Cells.Item(88, 9).Value = 9405.6  # I88
$ws.Cells.Item(88, 11).Value = 28216.8  # K88
$ws.Cells.Item(88, 13).Value = -27788.8  # M88
$ws.Cells.Item(91, 8).Value = 9117.777  # H91
$ws.Cells.Item(91, 9).Value = 9405.6  # I91
$ws.Cells.Item(91, 11).Value = 28216.8  # K91
$ws.Cells.Item(91, 13).Value = -26734.8  # M91
$ws.Cells.Item(107, 8).Value = 1740.8334  # H107
$ws.Cells.Item(107, 9).Value = 2599  # I107
$ws.Cells.Item(107, 10).Value = 1311.75  # J107
$ws.Cells.Item(107, 11).Value = 7797  # K107
$ws.Cells.Item(107, 12).Value = 3935.25  # L107
$ws.Cells.Item(107, 13).Value = -5877  # M107
$ws.Cells.Item(107, 14).Value = -7775.25  # N107
$ws.Cells.Item(121, 8).Value = 6947127  # H121
$ws.Cells.Item(121, 10).Value = 3112.95  # J121
$ws.Cells.Item(121, 12).Value = 9338.849999999999  # L121
$ws.Cells.Item(121, 14).Value = -11958.85  # N121
$ws.Cells.Item(139, 8).Value = 5525.1816  # H139
$ws.Cells.Item(139, 9).Value = 3412.4285  # I139
$ws.Cells.Item(139, 11).Value = 10237.2855  # K139
$ws.Cells.Item(139, 13).Value = -5097.2855  # M139

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 7191216  # H80
$ws.Cells.Item(80, 9).Value = 15337834  # I80
$ws.Cells.Item(80, 10).Value = 3023.5293  # J80
$ws.Cells.Item(80, 11).Value = 15337834  # K80
$ws.Cells.Item(80, 12).Value = 3023.5293  # L80
$ws.Cells.Item(80, 13).Value = -15336836  # M80
$ws.Cells.Item(80, 14).Value = -5019.5293  # N80
$ws.Cells.Item(83, 8).Value = 7191216  # H83
$ws.Cells.Item(83, 9).Value = 15337834  # I83
$ws.Cells.Item(83, 10).Value = 3023.5293  # J83
$ws.Cells.Item(83, 11).Value = 76689170  # K83
$ws.Cells.Item(83, 12).Value = 15117.6465  # L83
$ws.Cells.Item(83, 13).Value = -76684178  # M83
$ws.Cells.Item(83, 14).Value = -25101.6465  # N83
$ws.Cells.Item(102, 8).Value = 6128.357  # H102
$ws.Cells.Item(102, 9).Value = 6395.25  # I102
$ws.Cells.Item(102, 11).Value = 6395.25  # K102
$ws.Cells.Item(102, 13).Value = -4773.25  # M102
$ws.Cells.Item(110, 8).Value = 127998.4  # H110
$ws.Cells.Item(110, 10).Value = 127998.4  # J110
$ws.Cells.Item(110, 12).Value = 127998.4  # L110
$ws.Cells.Item(110, 14).Value = -136178.4  # N110
$ws.Cells.Item(113, 8).Value = 3362.4  # H113
$ws.Cells.Item(113, 9).Value = 2766  # I113
$ws.Cells.Item(113, 10).Value = 5002.5  # J113
$ws.Cells.Item(113, 11).Value = 2766  # K113
$ws.Cells.Item(113, 12).Value = 5002.5  # L113
$ws.Cells.Item(113, 13).Value = -596  # M113
$ws.Cells.Item(113, 14).Value = -9342.5  # N113
$ws.Cells.Item(126, 8).Value = 5905  # H126
$ws.Cells.Item(126, 9).Value = 6345.3076  # I126
$ws.Cells.Item(126, 10).Value = 5428  # J126
$ws.Cells.Item(126, 11).Value = 19035.9228  # K126
$ws.Cells.Item(126, 12).Value = 16284  # L126
$ws.Cells.Item(126, 13).Value = -16565.9228  # M126
$ws.Cells.Item(126, 14).Value = -21224  # N126
$ws.Cells.Item(132, 8).Value = 5819.6665  # H132
$ws.Cells.Item(132, 9).Value = 2318.6667  # I132
$ws.Cells.Item(132, 11).Value = 6956.000100000001  # K132
$ws.Cells.Item(132, 13).Value = -4426.000100000001  # M132
$ws.Cells.Item(138, 8).Value = 182496.75  # H138
$ws.Cells.Item(138, 10).Value = 182496.75  # J138
$ws.Cells.Item(138, 12).Value = 182496.75  # L138
$ws.Cells.Item(138, 14).Value = -192776.75  # N138

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3684.4  # H16
$ws.Cells.Item(16, 10).Value = 4660.3335  # J16
$ws.Cells.Item(16, 12).Value = 4660.3335  # L16
$ws.Cells.Item(16, 14).Value = -5000.3335  # N16
$ws.Cells.Item(30, 8).Value = 242.33333  # H30
$ws.Cells.Item(30, 9).Value = 242.33333  # I30
$ws.Cells.Item(30, 11).Value = 242.33333  # K30
$ws.Cells.Item(30, 13).Value = -134.33333  # M30
$ws.Cells.Item(40, 8).Value = 9132.223  # H40
$ws.Cells.Item(40, 9).Value = 8170.7144  # I40
$ws.Cells.Item(40, 11).Value = 8170.7144  # K40
$ws.Cells.Item(40, 13).Value = -8034.7144  # M40
$ws.Cells.Item(46, 8).Value = 1836.1428  # H46
$ws.Cells.Item(46, 9).Value = 725  # I46
$ws.Cells.Item(46, 10).Value = 2280.6  # J46
$ws.Cells.Item(46, 11).Value = 725  # K46
$ws.Cells.Item(46, 12).Value = 2280.6  # L46
$ws.Cells.Item(46, 13).Value = -537  # M46
$ws.Cells.Item(46, 14).Value = -2656.6  # N46
$ws.Cells.Item(68, 8).Value = 1798.5714  # H68
$ws.Cells.Item(68, 9).Value = 1800  # I68
$ws.Cells.Item(68, 11).Value = 1800  # K68
$ws.Cells.Item(68, 13).Value = -1051  # M68
$ws.Cells.Item(71, 8).Value = 1798.5714  # H71
$ws.Cells.Item(71, 9).Value = 1800  # I71
$ws.Cells.Item(71, 11).Value = 9000  # K71
$ws.Cells.Item(71, 13).Value = -5256  # M71
$ws.Cells.Item(128, 8).Value = 159999  # H128
$ws.Cells.Item(128, 10).Value = 159999  # J128
$ws.Cells.Item(128, 12).Value = 159999  # L128
$ws.Cells.Item(128, 14).Value = -169959  # N128
$ws.Cells.Item(132, 8).Value = 59475.844  # H132
$ws.Cells.Item(132, 9).Value = 69377.69  # I132
$ws.Cells.Item(132, 10).Value = 6666  # J132
$ws.Cells.Item(132, 11).Value = 208133.07  # K132
$ws.Cells.Item(132, 12).Value = 19998  # L132
$ws.Cells.Item(132, 13).Value = -205603.07  # M132
$ws.Cells.Item(132, 14).Value = -25058  # N132
$ws.Cells.Item(136, 8).Value = 6214913.5  # H136
$ws.Cells.Item(136, 10).Value = 13320.077  # J136
$ws.Cells.Item(136, 12).Value = 39960.231  # L136
$ws.Cells.Item(136, 14).Value = -45060.231  # N136
$ws.Cells.Item(137, 8).Value = 81139.125  # H137
$ws.Cells.Item(137, 10).Value = 87103.28999999999  # J137
$ws.Cells.Item(137, 12).Value = 87103.28999999999  # L137
$ws.Cells.Item(137, 14).Value = -97303.28999999999  # N137
$ws.Cells.Item(141, 8).Value = 81996.60000000001  # H141
$ws.Cells.Item(141, 10).Value = 81996.60000000001  # J141
$ws.Cells.Item(141, 12).Value = 81996.60000000001  # L141
$ws.Cells.Item(141, 14).Value = -92356.60000000001  # N141

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 9608  # H41
$ws.Cells.Item(41, 9).Value = 11162.333  # I41
$ws.Cells.Item(41, 10).Value = 8675.4  # J41
$ws.Cells.Item(41, 11).Value = 11162.333  # K41
$ws.Cells.Item(41, 12).Value = 8675.4  # L41
$ws.Cells.Item(41, 13).Value = -10772.333  # M41
$ws.Cells.Item(41, 14).Value = -9455.4  # N41
$ws.Cells.Item(46, 8).Value = 78824.27  # H46
$ws.Cells.Item(46, 10).Value = 78824.27  # J46
$ws.Cells.Item(46, 12).Value = 78824.27  # L46
$ws.Cells.Item(46, 14).Value = -79286.27  # N46
$ws.Cells.Item(81, 8).Value = 53153.05  # H81
$ws.Cells.Item(81, 9).Value = 2242.5454  # I81
$ws.Cells.Item(81, 10).Value = 115377  # J81
$ws.Cells.Item(81, 11).Value = 4485.0908  # K81
$ws.Cells.Item(81, 12).Value = 230754  # L81
$ws.Cells.Item(81, 13).Value = -3424.0908  # M81
$ws.Cells.Item(81, 14).Value = -232876  # N81
$ws.Cells.Item(84, 8).Value = 53153.05  # H84
$ws.Cells.Item(84, 9).Value = 2242.5454  # I84
$ws.Cells.Item(84, 10).Value = 115377  # J84
$ws.Cells.Item(84, 11).Value = 22425.454  # K84
$ws.Cells.Item(84, 12).Value = 1153770  # L84
$ws.Cells.Item(84, 13).Value = -17121.454  # M84
$ws.Cells.Item(84, 14).Value = -1164378  # N84
$ws.Cells.Item(104, 8).Value = 32911.6  # H104
$ws.Cells.Item(104, 10).Value = 32911.6  # J104
$ws.Cells.Item(104, 12).Value = 32911.6  # L104
$ws.Cells.Item(104, 14).Value = -39899.6  # N104
$ws.Cells.Item(107, 8).Value = 1462.3334  # H107
$ws.Cells.Item(107, 10).Value = 266.4  # J107
$ws.Cells.Item(107, 12).Value = 799.1999999999999  # L107
$ws.Cells.Item(107, 14).Value = -4639.2  # N107
$ws.Cells.Item(122, 8).Value = 14045  # H122
$ws.Cells.Item(122, 9).Value = 8309.6875  # I122
$ws.Cells.Item(122, 11).Value = 24929.0625  # K122
$ws.Cells.Item(122, 13).Value = -22479.0625  # M122
$ws.Cells.Item(134, 8).Value = 78824.27  # H134
$ws.Cells.Item(134, 10).Value = 78824.27  # J134
$ws.Cells.Item(134, 12).Value = 236472.81  # L134
$ws.Cells.Item(134, 14).Value = -241542.81  # N134
$ws.Cells.Item(136, 8).Value = 3357.4546  # H136
$ws.Cells.Item(136, 9).Value = 2048.5833  # I136
$ws.Cells.Item(136, 11).Value = 6145.749899999999  # K136
$ws.Cells.Item(136, 13).Value = -3595.749899999999  # M136
